$wb = $excel.ActiveWorkbook

# --- Sheet: Sim_5yr ---
$ws1 = $wb.Worksheets.Item("Sim_5yr")

$ws1.Range("B2").Value = -0.004282840336034826
$ws1.Range("C2").Value = 0.02302100293418834
$ws1.Range("D2").Value = 0.08579924713921762
$ws1.Range("E2").Value = 0.22516039235143

$ws1.Range("B3").Value = -0.01855104846907478
$ws1.Range("C3").Value = 0.005422169068455969
$ws1.Range("D3").Value = 0.07479120730010015
$ws1.Range("E3").Value = 0.2378289229639267

$ws1.Range("B4").Value = -0.02613736430319635
$ws1.Range("C4").Value = -0.007154220427063259
$ws1.Range("D4").Value = 0.05967910789027422
$ws1.Range("E4").Value = 0.2305681723758092

$ws1.Range("B5").Value = -0.01955700770304824
$ws1.Range("C5").Value = -0.001125629226679811
$ws1.Range("D5").Value = 0.06478498766028414
$ws1.Range("E5").Value = 0.2380158442442341

$ws1.Range("B6").Value = -0.00854105973134241
$ws1.Range("C6").Value = 0.01318274608469809
$ws1.Range("D6").Value = 0.08160447448843923
$ws1.Range("E6").Value = 0.255838728315701

$ws1.Range("B7").Value = -0.00006230550227516185
$ws1.Range("C7").Value = 0.02619313379369959
$ws1.Range("D7").Value = 0.09844080080017134
$ws1.Range("E7").Value = 0.2726153113653246

$ws1.Range("B8").Value = 0.004033704328844627
$ws1.Range("C8").Value = 0.03322454103344932
$ws1.Range("D8").Value = 0.1078608596092078
$ws1.Range("E8").Value = 0.2812938098090068

$ws1.Range("B9").Value = 0.006306681806613115
$ws1.Range("C9").Value = 0.03761272248300324
$ws1.Range("D9").Value = 0.1133327832638251
$ws1.Range("E9").Value = 0.2853922200211083

$ws1.Range("B10").Value = 0.001343701675643187
$ws1.Range("C10").Value = 0.02707462818091589
$ws1.Range("D10").Value = 0.09615231592399624
$ws1.Range("E10").Value = 0.2628644012402884

# --- Sheet: Sim_10yr ---
$ws2 = $wb.Worksheets.Item("Sim_10yr")

$ws2.Range("B2").Value = -0.004282840336034826
$ws2.Range("C2").Value = 0.02302100293418834
$ws2.Range("D2").Value = 0.08579924713921762
$ws2.Range("E2").Value = 0.22516039235143

$ws2.Range("B3").Value = -0.0226730292340869
$ws2.Range("C3").Value = -0.001373785615127518
$ws2.Range("D3").Value = 0.0659747206858326
$ws2.Range("E3").Value = 0.2316682306050966

$ws2.Range("B4").Value = -0.01353279462471494
$ws2.Range("C4").Value = 0.006319617599752947
$ws2.Range("D4").Value = 0.07269332730951769
$ws2.Range("E4").Value = 0.2450556704557504

$ws2.Range("B5").Value = 0.001959036517865503
$ws2.Range("C5").Value = 0.02945572141998663
$ws2.Range("D5").Value = 0.1025130401938906
$ws2.Range("E5").Value = 0.2757024079197318

$ws2.Range("B6").Value = 0.003841510628239411
$ws2.Range("C6").Value = 0.03250357040539691
$ws2.Range("D6").Value = 0.1051671806429676
$ws2.Range("E6").Value = 0.2748122718004537
